$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.3407354241124098; C = 0.5844350863438295; D = 0.4442199068852112; E = 0.6664982422221466; F = 0.5885134887136687; G = 19 }
    3  = @{ B = 0.1726228500093137; C = 0.4551251828485784; D = 0.3113871153733175; E = 0.5580207123156966; F = 0.5460334188492619; G = 18 }
    4  = @{ B = 0.2034452847961829; C = 0.3989700209739589; D = 0.2615574828328119; E = 0.5114269085928232; F = 0.483661004847352;  G = 17 }
    5  = @{ B = 0.3635533387143693; C = 0.3920234804570175; D = 0.2041752638615354; E = 0.4518575703266854; F = 0.2771362769876616; G = 16 }
    6  = @{ B = 0.3372887961726861; C = 0.367293011109447;  D = 0.1918266328756511; E = 0.4379801740668761; F = 0.289203773041539;  G = 15 }
    7  = @{ B = 0.3465491651461942; C = 0.3849987877080897; D = 0.2040158891459988; E = 0.4516811808632266; F = 0.3006242113637049; G = 14 }
    8  = @{ B = 0.3728115640212805; C = 0.4172293592727421; D = 0.2174429248205509; E = 0.4663077576242443; F = 0.2915344481050866; G = 13 }
    9  = @{ B = 0.4140746927647911; C = 0.4327887006433242; D = 0.2237006186653438; E = 0.4729699976376343; F = 0.2387302033153843; G = 12 }
    10 = @{ B = 0.3853625486457603; C = 0.4009980026591189; D = 0.1935302347147137; E = 0.4399207141232539; F = 0.2225500727871226; G = 11 }
    11 = @{ B = 0.3626069129865366; C = 0.3740635702753569; D = 0.1753426339092577; E = 0.4187393388604153; F = 0.2207534083381827; G = 10 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Range("B$row").Value = $rowData.B
    $ws.Range("C$row").Value = $rowData.C
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("E$row").Value = $rowData.E
    $ws.Range("F$row").Value = $rowData.F
    $ws.Range("G$row").Value = $rowData.G
}

$wb.Save()
